$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) and Volume(1h) (E) columns to Text format before writing,
# so numeric-looking strings (e.g. "19.63") are not auto-converted to numbers.
$colD = $ws.Range("D2:D51")
$colE = $ws.Range("E2:E51")
$colD.NumberFormat = "@"
$colE.NumberFormat = "@"

$ws.Range("D2").Value = "26.664.29"
$ws.Range("E2").Value = "  -1.53%  "
$ws.Range("D3").Value = "1.594.09"
$ws.Range("E3").Value = "  -1.69%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("E5").Value = "  -1.38%  "
$ws.Range("D6").Value = "0.509"
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  -1.60%  "
$ws.Range("E9").Value = "  -1.58%  "
$ws.Range("D10").Value = "19.63"
$ws.Range("E10").Value = "  -1.37%  "
$ws.Range("D11").Value = "0.0836"
$ws.Range("E11").Value = "  -0.65%  "
$ws.Range("D12").Value = "1.818.58"
$ws.Range("E12").Value = "  -1.62%  "
$ws.Range("D13").Value = "1.600.25"
$ws.Range("E13").Value = "  -1.32%  "
$ws.Range("E14").Value = "  -2.32%  "
$ws.Range("E15").Value = "  -2.89%  "
$ws.Range("D16").Value = "64.74"
$ws.Range("E16").Value = "  +0.42%  "
$ws.Range("D17").Value = "26.644.07"
$ws.Range("E17").Value = "  -1.46%  "
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("D19").Value = "208.33"
$ws.Range("E19").Value = "  -2.74%  "
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").Value = "6.71"
$ws.Range("E21").Value = "  -1.85%  "
$ws.Range("E22").Value = "  -2.32%  "
$ws.Range("E23").Value = "  +0.58%  "
$ws.Range("D24").Value = "8.89"
$ws.Range("E24").Value = "  -1.30%  "
$ws.Range("D25").Value = "146.51"
$ws.Range("E25").Value = "  -0.86%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "7.13"
$ws.Range("E27").Value = "  -3.63%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("D29").Value = "15.30"
$ws.Range("E29").Value = "  -1.26%  "
$ws.Range("D30").Value = "0.0502"
$ws.Range("E30").Value = "  -1.28%  "
$ws.Range("D31").Value = "1.16"
$ws.Range("E31").Value = "  -1.09%  "
$ws.Range("E32").Value = "  -2.99%  "
$ws.Range("D33").Value = "0.674"
$ws.Range("E33").Value = "  -4.22%  "
$ws.Range("E34").Value = "  -2.44%  "
$ws.Range("D35").Value = "1.299.29"
$ws.Range("E35").Value = "  -3.22%  "
$ws.Range("E36").Value = "  -0.53%  "
$ws.Range("E37").Value = "  -4.98%  "
$ws.Range("E38").Value = "  -2.76%  "
$ws.Range("D39").Value = "0.835"
$ws.Range("E39").Value = "  -0.68%  "
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "5.35"
$ws.Range("E42").Value = "  +0.56%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").Value = "2.19"
$ws.Range("E43").Value = "  -1.77%  "
$ws.Range("E44").Value = "  -0.34%  "
$ws.Range("D45").Value = "1.730.25"
$ws.Range("E45").Value = "  -1.67%  "
$ws.Range("D46").Value = "0.910"
$ws.Range("E46").Value = "  +6.82%  "
$ws.Range("D47").Value = "89.73"
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("D48").Value = "1.64"
$ws.Range("E48").Value = "  -0.37%  "
$ws.Range("D49").Value = "0.0984"
$ws.Range("E49").Value = "  -1.86%  "
$ws.Range("E50").Value = "  -1.71%  "
$ws.Range("D51").Value = "7.52"
$ws.Range("E51").Value = "  -0.74%  "

# Restore the default cell style so no formatting diff is introduced.
$colD.Style = "Normal"
$colE.Style = "Normal"
